$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'bev0F86L'
$ws.Cells.Item(2, 2).Value = '20/11/2024'
$ws.Cells.Item(2, 3).Value = '21:30'
$ws.Cells.Item(2, 4).Value = 'ARGENTINA - TORNEO BETANO'
$ws.Cells.Item(2, 5).Value = 'Boca Juniors'
$ws.Cells.Item(2, 6).Value = 'Union de Santa Fe'
$ws.Cells.Item(2, 7).Value = 2.05
$ws.Cells.Item(2, 8).Value = 3.2
$ws.Cells.Item(2, 9).Value = 3.9
$ws.Cells.Item(2, 10).Value = 2.88
$ws.Cells.Item(2, 11).Value = 1.95
$ws.Cells.Item(2, 12).Value = 4.75
$ws.Cells.Item(2, 13).Value = 1.1
$ws.Cells.Item(2, 14).Value = 7
$ws.Cells.Item(2, 15).Value = 1.5
$ws.Cells.Item(2, 16).Value = 2.5
$ws.Cells.Item(2, 17).Value = 2.5
$ws.Cells.Item(2, 18).Value = 1.5
$ws.Cells.Item(2, 19).Value = 1.57
$ws.Cells.Item(2, 20).Value = 2.25
$ws.Cells.Item(2, 21).Value = 2.2
$ws.Cells.Item(2, 22).Value = 1.62
$ws.Cells.Item(2, 23).Value = 5.5
$ws.Cells.Item(2, 24).Value = 8.5
$ws.Cells.Item(2, 25).Value = 9.5
$ws.Cells.Item(2, 26).Value = 17
$ws.Cells.Item(2, 27).Value = 21
$ws.Cells.Item(2, 28).Value = 41
$ws.Cells.Item(2, 29).Value = 6.5
$ws.Cells.Item(2, 30).Value = 6.5
$ws.Cells.Item(2, 31).Value = 19
$ws.Cells.Item(2, 32).Value = 81
$ws.Cells.Item(2, 33).Value = 501
$ws.Cells.Item(2, 34).Value = 8.5
$ws.Cells.Item(2, 35).Value = 19
$ws.Cells.Item(2, 36).Value = 15
$ws.Cells.Item(2, 37).Value = 41
$ws.Cells.Item(2, 38).Value = 41
$ws.Cells.Item(2, 39).Value = 51
$ws.Cells.Item(2, 40).Value = 3.75
$ws.Cells.Item(2, 41).Value = 12
$ws.Cells.Item(2, 42).Value = 29
$ws.Cells.Item(2, 43).Value = 41
$ws.Cells.Item(2, 44).Value = 81
$ws.Cells.Item(2, 45).Value = 251
$ws.Cells.Item(2, 46).Value = 2.25
$ws.Cells.Item(2, 47).Value = 9.5
$ws.Cells.Item(2, 48).Value = 81
$ws.Cells.Item(2, 49).Value = 5.5
$ws.Cells.Item(2, 50).Value = 23
$ws.Cells.Item(2, 51).Value = 41
$ws.Cells.Item(2, 52).Value = 81
$ws.Cells.Item(2, 53).Value = 126
$ws.Cells.Item(2, 54).Value = 401
$ws.Cells.Item(2, 55).Value = 126
$ws.Cells.Item(2, 56).Value = 126

# Row 3
$ws.Cells.Item(3, 1).Value = 'OAukwB1L'
$ws.Cells.Item(3, 2).Value = '20/11/2024'
$ws.Cells.Item(3, 3).Value = '21:00'
$ws.Cells.Item(3, 4).Value = 'BOLIVIA - DIVISION PROFESIONAL'
$ws.Cells.Item(3, 5).Value = 'Universitario de Vinto'
$ws.Cells.Item(3, 6).Value = 'Oriente Petrolero'
$ws.Cells.Item(3, 7).Value = 2.15
$ws.Cells.Item(3, 8).Value = 3.3
$ws.Cells.Item(3, 9).Value = 3.3
$ws.Cells.Item(3, 10).Value = 2.88
$ws.Cells.Item(3, 11).Value = 2.05
$ws.Cells.Item(3, 12).Value = 4
$ws.Cells.Item(3, 13).Value = 1.06
$ws.Cells.Item(3, 14).Value = 10
$ws.Cells.Item(3, 15).Value = 1.33
$ws.Cells.Item(3, 16).Value = 3.25
$ws.Cells.Item(3, 17).Value = 2.1
$ws.Cells.Item(3, 18).Value = 1.7
$ws.Cells.Item(3, 19).Value = 1.5
$ws.Cells.Item(3, 20).Value = 2.5
$ws.Cells.Item(3, 21).Value = 1.91
$ws.Cells.Item(3, 22).Value = 1.8
$ws.Cells.Item(3, 23).Value = 6.5
$ws.Cells.Item(3, 24).Value = 9.5
$ws.Cells.Item(3, 25).Value = 9.5
$ws.Cells.Item(3, 26).Value = 19
$ws.Cells.Item(3, 27).Value = 19
$ws.Cells.Item(3, 28).Value = 34
$ws.Cells.Item(3, 29).Value = 8.5
$ws.Cells.Item(3, 30).Value = 6.5
$ws.Cells.Item(3, 31).Value = 17
$ws.Cells.Item(3, 32).Value = 51
$ws.Cells.Item(3, 33).Value = 351
$ws.Cells.Item(3, 34).Value = 9
$ws.Cells.Item(3, 35).Value = 15
$ws.Cells.Item(3, 36).Value = 12
$ws.Cells.Item(3, 37).Value = 34
$ws.Cells.Item(3, 38).Value = 29
$ws.Cells.Item(3, 39).Value = 41
$ws.Cells.Item(3, 40).Value = 4
$ws.Cells.Item(3, 41).Value = 12
$ws.Cells.Item(3, 42).Value = 23
$ws.Cells.Item(3, 43).Value = 41
$ws.Cells.Item(3, 44).Value = 67
$ws.Cells.Item(3, 45).Value = 201
$ws.Cells.Item(3, 46).Value = 2.5
$ws.Cells.Item(3, 47).Value = 8.5
$ws.Cells.Item(3, 48).Value = 67
$ws.Cells.Item(3, 49).Value = 5
$ws.Cells.Item(3, 50).Value = 21
$ws.Cells.Item(3, 51).Value = 29
$ws.Cells.Item(3, 52).Value = 67
$ws.Cells.Item(3, 53).Value = 101
$ws.Cells.Item(3, 54).Value = 251
$ws.Cells.Item(3, 55).Value = ""
$ws.Cells.Item(3, 56).Value = ""

# Row 4
$ws.Cells.Item(4, 1).Value = 'l8uhrXe8'
$ws.Cells.Item(4, 2).Value = '20/11/2024'
$ws.Cells.Item(4, 3).Value = '21:30'
$ws.Cells.Item(4, 4).Value = 'BRAZIL - SERIE A BETANO'
$ws.Cells.Item(4, 5).Value = 'Atletico-MG'
$ws.Cells.Item(4, 6).Value = 'Botafogo RJ'
$ws.Cells.Item(4, 7).Value = 3.5
$ws.Cells.Item(4, 8).Value = 3.2
$ws.Cells.Item(4, 9).Value = 2.2
$ws.Cells.Item(4, 10).Value = 4
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 3
$ws.Cells.Item(4, 13).Value = 1.08
$ws.Cells.Item(4, 14).Value = 8
$ws.Cells.Item(4, 15).Value = 1.4
$ws.Cells.Item(4, 16).Value = 3
$ws.Cells.Item(4, 17).Value = 2.3
$ws.Cells.Item(4, 18).Value = 1.6
$ws.Cells.Item(4, 19).Value = 1.5
$ws.Cells.Item(4, 20).Value = 2.5
$ws.Cells.Item(4, 21).Value = 1.95
$ws.Cells.Item(4, 22).Value = 1.8
$ws.Cells.Item(4, 23).Value = 8.5
$ws.Cells.Item(4, 24).Value = 17
$ws.Cells.Item(4, 25).Value = 13
$ws.Cells.Item(4, 26).Value = 41
$ws.Cells.Item(4, 27).Value = 29
$ws.Cells.Item(4, 28).Value = 41
$ws.Cells.Item(4, 29).Value = 7.5
$ws.Cells.Item(4, 30).Value = 6
$ws.Cells.Item(4, 31).Value = 17
$ws.Cells.Item(4, 32).Value = 51
$ws.Cells.Item(4, 33).Value = 401
$ws.Cells.Item(4, 34).Value = 6.5
$ws.Cells.Item(4, 35).Value = 9.5
$ws.Cells.Item(4, 36).Value = 9.5
$ws.Cells.Item(4, 37).Value = 21
$ws.Cells.Item(4, 38).Value = 21
$ws.Cells.Item(4, 39).Value = 34
$ws.Cells.Item(4, 40).Value = 5
$ws.Cells.Item(4, 41).Value = 21
$ws.Cells.Item(4, 42).Value = 34
$ws.Cells.Item(4, 43).Value = 67
$ws.Cells.Item(4, 44).Value = 101
$ws.Cells.Item(4, 45).Value = 251
$ws.Cells.Item(4, 46).Value = 2.5
$ws.Cells.Item(4, 47).Value = 8.5
$ws.Cells.Item(4, 48).Value = 67
$ws.Cells.Item(4, 49).Value = 4
$ws.Cells.Item(4, 50).Value = 13
$ws.Cells.Item(4, 51).Value = 26
$ws.Cells.Item(4, 52).Value = 41
$ws.Cells.Item(4, 53).Value = 67
$ws.Cells.Item(4, 54).Value = 201
$ws.Cells.Item(4, 55).Value = 126
$ws.Cells.Item(4, 56).Value = 126

# Row 5
$ws.Cells.Item(5, 1).Value = 'IZ3qIEYa'
$ws.Cells.Item(5, 2).Value = '20/11/2024'
$ws.Cells.Item(5, 3).Value = '20:30'
$ws.Cells.Item(5, 4).Value = 'COLOMBIA - PRIMERA A'
$ws.Cells.Item(5, 5).Value = 'Millonarios'
$ws.Cells.Item(5, 6).Value = 'Dep. Pasto'
$ws.Cells.Item(5, 7).Value = 1.73
$ws.Cells.Item(5, 8).Value = 3.25
$ws.Cells.Item(5, 9).Value = 5.75
$ws.Cells.Item(5, 10).Value = 2.5
$ws.Cells.Item(5, 11).Value = 1.95
$ws.Cells.Item(5, 12).Value = 6
$ws.Cells.Item(5, 13).Value = 1.11
$ws.Cells.Item(5, 14).Value = 6.5
$ws.Cells.Item(5, 15).Value = 1.5
$ws.Cells.Item(5, 16).Value = 2.5
$ws.Cells.Item(5, 17).Value = 2.6
$ws.Cells.Item(5, 18).Value = 1.48
$ws.Cells.Item(5, 19).Value = 1.57
$ws.Cells.Item(5, 20).Value = 2.25
$ws.Cells.Item(5, 21).Value = 2.5
$ws.Cells.Item(5, 22).Value = 1.5
$ws.Cells.Item(5, 23).Value = 5
$ws.Cells.Item(5, 24).Value = 6.5
$ws.Cells.Item(5, 25).Value = 9.5
$ws.Cells.Item(5, 26).Value = 13
$ws.Cells.Item(5, 27).Value = 19
$ws.Cells.Item(5, 28).Value = 41
$ws.Cells.Item(5, 29).Value = 6
$ws.Cells.Item(5, 30).Value = 6.5
$ws.Cells.Item(5, 31).Value = 23
$ws.Cells.Item(5, 32).Value = 101
$ws.Cells.Item(5, 33).Value = 201
$ws.Cells.Item(5, 34).Value = 11
$ws.Cells.Item(5, 35).Value = 26
$ws.Cells.Item(5, 36).Value = 21
$ws.Cells.Item(5, 37).Value = 67
$ws.Cells.Item(5, 38).Value = 51
$ws.Cells.Item(5, 39).Value = 67
$ws.Cells.Item(5, 40).Value = 3.4
$ws.Cells.Item(5, 41).Value = 9.5
$ws.Cells.Item(5, 42).Value = 29
$ws.Cells.Item(5, 43).Value = 34
$ws.Cells.Item(5, 44).Value = 67
$ws.Cells.Item(5, 45).Value = 301
$ws.Cells.Item(5, 46).Value = 2.25
$ws.Cells.Item(5, 47).Value = 10
$ws.Cells.Item(5, 48).Value = 81
$ws.Cells.Item(5, 49).Value = 7
$ws.Cells.Item(5, 50).Value = 34
$ws.Cells.Item(5, 51).Value = 41
$ws.Cells.Item(5, 52).Value = 151
$ws.Cells.Item(5, 53).Value = 201
$ws.Cells.Item(5, 54).Value = 501
$ws.Cells.Item(5, 55).Value = 126
$ws.Cells.Item(5, 56).Value = 126

# Row 6
$ws.Cells.Item(6, 1).Value = 'E1chGh3C'
$ws.Cells.Item(6, 2).Value = '20/11/2024'
$ws.Cells.Item(6, 3).Value = '22:30'
$ws.Cells.Item(6, 4).Value = 'COLOMBIA - PRIMERA A'
$ws.Cells.Item(6, 5).Value = 'Atl. Nacional'
$ws.Cells.Item(6, 6).Value = 'Santa Fe'
$ws.Cells.Item(6, 7).Value = 2
$ws.Cells.Item(6, 8).Value = 3
$ws.Cells.Item(6, 9).Value = 4.33
$ws.Cells.Item(6, 10).Value = 2.75
$ws.Cells.Item(6, 11).Value = 1.95
$ws.Cells.Item(6, 12).Value = 4.75
$ws.Cells.Item(6, 13).Value = 1.1
$ws.Cells.Item(6, 14).Value = 7
$ws.Cells.Item(6, 15).Value = 1.44
$ws.Cells.Item(6, 16).Value = 2.63
$ws.Cells.Item(6, 17).Value = 2.4
$ws.Cells.Item(6, 18).Value = 1.53
$ws.Cells.Item(6, 19).Value = 1.53
$ws.Cells.Item(6, 20).Value = 2.38
$ws.Cells.Item(6, 21).Value = 2.1
$ws.Cells.Item(6, 22).Value = 1.67
$ws.Cells.Item(6, 23).Value = 6
$ws.Cells.Item(6, 24).Value = 8.5
$ws.Cells.Item(6, 25).Value = 9.5
$ws.Cells.Item(6, 26).Value = 17
$ws.Cells.Item(6, 27).Value = 19
$ws.Cells.Item(6, 28).Value = 34
$ws.Cells.Item(6, 29).Value = 7
$ws.Cells.Item(6, 30).Value = 6.5
$ws.Cells.Item(6, 31).Value = 19
$ws.Cells.Item(6, 32).Value = 67
$ws.Cells.Item(6, 33).Value = 1250
$ws.Cells.Item(6, 34).Value = 9
$ws.Cells.Item(6, 35).Value = 19
$ws.Cells.Item(6, 36).Value = 15
$ws.Cells.Item(6, 37).Value = 41
$ws.Cells.Item(6, 38).Value = 41
$ws.Cells.Item(6, 39).Value = 41
$ws.Cells.Item(6, 40).Value = 3.75
$ws.Cells.Item(6, 41).Value = 12
$ws.Cells.Item(6, 42).Value = 26
$ws.Cells.Item(6, 43).Value = 41
$ws.Cells.Item(6, 44).Value = 67
$ws.Cells.Item(6, 45).Value = 251
$ws.Cells.Item(6, 46).Value = 2.38
$ws.Cells.Item(6, 47).Value = 9
$ws.Cells.Item(6, 48).Value = 67
$ws.Cells.Item(6, 49).Value = 5.5
$ws.Cells.Item(6, 50).Value = 23
$ws.Cells.Item(6, 51).Value = 34
$ws.Cells.Item(6, 52).Value = 81
$ws.Cells.Item(6, 53).Value = 126
$ws.Cells.Item(6, 54).Value = 351
$ws.Cells.Item(6, 55).Value = 126
$ws.Cells.Item(6, 56).Value = 126

# Row 7
$ws.Cells.Item(7, 1).Value = 'MREUeAWQ'
$ws.Cells.Item(7, 2).Value = '20/11/2024'
$ws.Cells.Item(7, 3).Value = '22:00'
$ws.Cells.Item(7, 4).Value = 'MEXICO - LIGA DE EXPANSION MX'
$ws.Cells.Item(7, 5).Value = 'Celaya'
$ws.Cells.Item(7, 6).Value = 'Tapatio'
$ws.Cells.Item(7, 7).Value = 1.6
$ws.Cells.Item(7, 8).Value = 3.75
$ws.Cells.Item(7, 9).Value = 5.1
$ws.Cells.Item(7, 10).Value = 2.15
$ws.Cells.Item(7, 11).Value = 2.18
$ws.Cells.Item(7, 12).Value = 5.1
$ws.Cells.Item(7, 13).Value = 1.02
$ws.Cells.Item(7, 14).Value = 10
$ws.Cells.Item(7, 15).Value = 1.26
$ws.Cells.Item(7, 16).Value = 3.15
$ws.Cells.Item(7, 17).Value = 1.83
$ws.Cells.Item(7, 18).Value = 1.88
$ws.Cells.Item(7, 19).Value = 1.38
$ws.Cells.Item(7, 20).Value = 2.57
$ws.Cells.Item(7, 21).Value = 1.8
$ws.Cells.Item(7, 22).Value = 1.8
$ws.Cells.Item(7, 23).Value = 6.7
$ws.Cells.Item(7, 24).Value = 7.3
$ws.Cells.Item(7, 25).Value = 8
$ws.Cells.Item(7, 26).Value = 11.5
$ws.Cells.Item(7, 27).Value = 13
$ws.Cells.Item(7, 28).Value = 27
$ws.Cells.Item(7, 29).Value = 10.5
$ws.Cells.Item(7, 30).Value = 7.3
$ws.Cells.Item(7, 31).Value = 16.5
$ws.Cells.Item(7, 32).Value = 75
$ws.Cells.Item(7, 33).Value = 600
$ws.Cells.Item(7, 34).Value = 14
$ws.Cells.Item(7, 35).Value = 30
$ws.Cells.Item(7, 36).Value = 16.5
$ws.Cells.Item(7, 37).Value = 100
$ws.Cells.Item(7, 38).Value = 50
$ws.Cells.Item(7, 39).Value = 50
$ws.Cells.Item(7, 40).Value = 3.4
$ws.Cells.Item(7, 41).Value = 7.7
$ws.Cells.Item(7, 42).Value = 17.5
$ws.Cells.Item(7, 43).Value = 25
$ws.Cells.Item(7, 44).Value = 55
$ws.Cells.Item(7, 45).Value = 250
$ws.Cells.Item(7, 46).Value = 2.6
$ws.Cells.Item(7, 47).Value = 7.6
$ws.Cells.Item(7, 48).Value = 70
$ws.Cells.Item(7, 49).Value = 6.6
$ws.Cells.Item(7, 50).Value = 28
$ws.Cells.Item(7, 51).Value = 32
$ws.Cells.Item(7, 52).Value = 175
$ws.Cells.Item(7, 53).Value = 200
$ws.Cells.Item(7, 54).Value = 450
$ws.Cells.Item(7, 55).Value = 500
$ws.Cells.Item(7, 56).Value = 51

# Remove the old row 8 (data now consolidated into rows 2-7)
$ws.Rows.Item(8).Delete()

"Done. UsedRange: " + $ws.UsedRange.Address()